$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "26.860.40"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "1.869.18"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'305.01"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'0.5073"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("D8").Value = "'0.3656"
$ws.Range("E8").Value = "  -2.77%  "
$ws.Range("D9").Value = "'0.07191"
$ws.Range("E9").Value = "  +0.47%  "
$ws.Range("D10").Value = "'0.8926"
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "'0.07523"
$ws.Range("E12").Value = "  -1.06%  "
$ws.Range("D13").Value = "1.877.90"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "'95.02"
$ws.Range("E14").Value = "  +6.19%  "
$ws.Range("D15").Value = "'5.225"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "'0.000008496"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  +0.62%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "26.926.01"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").Value = "2.114.37"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("D24").Value = "'6.381"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "'2.101"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'113.35"
$ws.Range("E29").Value = "  +0.57%  "
$ws.Range("D30").Value = "'4.701"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("D31").Value = "'4.727"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "'0.09141"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -1.16%  "
$ws.Range("D34").Value = "'0.7468"
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("D35").Value = "'2.983"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "'3.230"
$ws.Range("E37").Value = "  +5.72%  "
$ws.Range("D38").Value = "'2.529"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("E39").Value = "  +4.80%  "
$ws.Range("E40").Value = "  -2.28%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "'6.610"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("D43").Value = "'115.48"
$ws.Range("E43").Value = "  -0.14%  "
$ws.Range("D44").Value = "'8.590"
$ws.Range("E44").Value = "  +3.57%  "
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").Value = "'0.4745"
$ws.Range("E46").Value = "  +2.20%  "
$ws.Range("D47").Value = "'1.0000"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").Value = "'10.10"
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "'36.91"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("D51").Value = "'63.03"
$ws.Range("E51").Value = "  -1.06%  "
